# Add a new slide (slide 2) using the same "Title, Content" layout as slide 1.
$p = $ppt.ActivePresentation

$design = $p.Designs.Item(1)
$customLayout = $design.SlideMaster.CustomLayouts.Item(3)   # "Title, Content"
$s = $p.Slides.AddSlide(2, $customLayout)

# ---------------------------------------------------------------------------
# Body / text placeholder (shape 1)
# ---------------------------------------------------------------------------
$body = $s.Shapes.Item(1)

$body.Left = 39.68503937007874
$body.Top = 82.16795275590552
$body.Width = 714.3023622047244
$body.Height = 487.6131496062992

$bodyTf = $body.TextFrame
$bodyTf.AutoSize = 0   # ppAutoSizeNone -> <a:noAutofit/>

$lines = @(
    "",
    "Project management processes fall into five groups:",
    "",
    "Initiating",
    "Planning",
    "Executing",
    "Monitoring and Controlling",
    "Closing",
    "",
    "Project management knowledge draws on ten areas:",
    "",
    "Integration",
    "Scope",
    "Time",
    "Cost",
    "Quality",
    "Procurement",
    "Human resources",
    "Communications",
    "Risk management",
    "Stakeholder management",
    ""
)

$bodyTr = $bodyTf.TextRange
$bodyTr.Text = [string]::Join([char]13, $lines)

# Base font formatting for every paragraph in the body.
$bodyTr.Font.Name = "Times New Roman"
$bodyTr.Font.NameComplexScript = "Times New Roman"
$bodyTr.Font.Size = 17

# Bold header lines.
$boldParagraphs = @(2, 10)
foreach ($i in $boldParagraphs) {
    $bodyTr.Paragraphs($i, 1).Font.Bold = $true
}

# Bulleted lines get a bullet character + indent level.
$bulletParagraphs = @(4, 5, 6, 7, 8, 9, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21)
foreach ($i in $bulletParagraphs) {
    $para = $bodyTr.Paragraphs($i, 1)
    $para.IndentLevel = 1
    $bullet = $para.ParagraphFormat.Bullet
    $bullet.Visible = $true
    $bullet.Character = 8226
    $bullet.Font.Name = "Arial"
}

# ---------------------------------------------------------------------------
# Title placeholder (shape 2)
# ---------------------------------------------------------------------------
$title = $s.Shapes.Item(2)

$title.Left = 39.68503937007874
$title.Top = 23.725984251968505
$title.Width = 714.3023622047244
$title.Height = 41.43220472440945

$titleTr = $title.TextFrame.TextRange
$titleTr.Text = "Project Management"
$titleTr.Font.Size = 32
$titleTr.Font.Bold = $true
$titleTr.Font.Underline = $true
$titleTr.Font.Name = "Times New Roman"
$titleTr.Font.NameComplexScript = "Times New Roman"
$titleTr.Font.Color.RGB = 0x990000
$titleTr.ParagraphFormat.Alignment = 2  # ppAlignCenter
